$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)
$run = $para.Runs(1, 1)
$run.Text = "Philippe.Renevier-Gonin@ac-grenoble.fr"
